# Update cube metadata: Package4
#
# 1) Rename the worksheet to match the new report code.
# 2) Resize the data columns to the new layout (15 custom columns instead
#    of 19) and reset the now-unused trailing columns back to the
#    (approximate) default width.
# 3) Re-create the trailing blank rows that are now part of the sheet.
# 4) Move the selection to B19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename sheet -------------------------------------------------
$ws.Name = "Informe-01-010057-A-TC-TM-TP"

# --- 2. Column widths --------------------------------------------------
# Excel's ColumnWidth property is expressed in "characters"; the value
# actually persisted in the XML adds a fixed 5/6 character padding on
# top of the number of characters that fit (rounded to whole pixels).
# Offsetting our desired stored width by that padding reproduces the
# target widths as closely as this engine's pixel rounding allows.
$offset = 5.0 / 6.0

$widths = @{
    1  = 27.69
    2  = 44.5
    3  = 18.66
    4  = 55.2
    5  = 34.64
    6  = 36.31
    7  = 47.28
    8  = 27.69
    9  = 27.69
    10 = 27.69
    11 = 15.46
    12 = 46.44
    13 = 19.19
    14 = 20.05
    15 = 29.5
    16 = 11.52
    17 = 11.52
    18 = 11.52
    19 = 11.52
}

foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col] - $offset
}

# --- 3. Extra trailing blank rows --------------------------------------
for ($r = 7; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}

# --- 4. Selection --------------------------------------------------
$null = $ws.Range("B19").Select()
